# Corrige os valores de "indirect" na planilha (tabela de meta-analise em rede)
# conforme commit: "corrigi o valor de indirect na planilha cinema H -> L"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 2
$ws.Range("A2").Value = "imipramina"
$ws.Range("B2").Value = "-0.26 (-3.31;  2.80)"
$ws.Range("C2").Value = "."
$ws.Range("E2").Value = "2.40 (-0.76;  5.57)"
$ws.Range("F2").Value = "."
$ws.Range("H2").Value = "3.72 ( 0.40;  7.03)"
$ws.Range("I2").Value = "2.03 ( 1.32;  2.74)"

# Linha 3
$ws.Range("A3").Value = "0.00 (-1.16;  1.17)"
$ws.Range("B3").Value = "desipramina"
$ws.Range("C3").Value = "."
$ws.Range("E3").Value = "1.72 (-1.40;  4.83)"
$ws.Range("F3").Value = "."
$ws.Range("H3").Value = "2.67 (-0.52;  5.86)"
$ws.Range("I3").Value = "2.04 ( 1.05;  3.02)"

# Linha 4
$ws.Range("A4").Value = "0.56 (-0.83;  1.94)"
$ws.Range("B4").Value = "0.55 (-0.98;  2.09)"
$ws.Range("C4").Value = "mianserina"
$ws.Range("F4").Value = "."
$ws.Range("I4").Value = "1.47 ( 0.27;  2.66)"

# Linha 5
$ws.Range("A5").Value = "0.79 (-0.37;  1.95)"
$ws.Range("B5").Value = "0.79 (-0.54;  2.12)"
$ws.Range("C5").Value = "0.24 (-1.27;  1.74)"
$ws.Range("D5").Value = "fluoxetina"
$ws.Range("I5").Value = "1.23 ( 0.31;  2.15)"

# Linha 6
$ws.Range("A6").Value = "0.85 (-0.68;  2.39)"
$ws.Range("B6").Value = "0.85 (-0.79;  2.50)"
$ws.Range("C6").Value = "0.30 (-1.54;  2.13)"
$ws.Range("D6").Value = "0.06 (-1.61;  1.73)"
$ws.Range("E6").Value = "amitriptilina"
$ws.Range("H6").Value = "1.26 (-1.83;  4.34)"
$ws.Range("I6").Value = "1.34 (-0.13;  2.81)"

# Linha 7
$ws.Range("A7").Value = "1.63 (-0.55;  3.81)"
$ws.Range("B7").Value = "1.63 (-0.64;  3.90)"
$ws.Range("C7").Value = "1.08 (-1.31;  3.46)"
$ws.Range("D7").Value = "0.84 (-1.42;  3.10)"
$ws.Range("E7").Value = "0.78 (-1.71;  3.27)"
$ws.Range("F7").Value = "citalopram"
$ws.Range("I7").Value = "0.39 (-1.67;  2.45)"

# Linha 8
$ws.Range("A8").Value = "1.67 (-0.16;  3.49)"
$ws.Range("B8").Value = "1.66 (-0.27;  3.60)"
$ws.Range("C8").Value = "1.11 (-0.95;  3.18)"
$ws.Range("D8").Value = "0.88 (-1.04;  2.80)"
$ws.Range("E8").Value = "0.81 (-1.37;  3.00)"
$ws.Range("F8").Value = "0.04 (-2.63;  2.70)"
$ws.Range("G8").Value = "fluvoxamina"
$ws.Range("I8").Value = "0.35 (-1.33;  2.04)"

# Linha 9
$ws.Range("A9").Value = "1.88 ( 0.32;  3.43)"
$ws.Range("B9").Value = "1.87 ( 0.22;  3.53)"
$ws.Range("C9").Value = "1.32 (-0.52;  3.16)"
$ws.Range("D9").Value = "1.09 (-0.59;  2.76)"
$ws.Range("E9").Value = "1.02 (-0.86;  2.91)"
$ws.Range("F9").Value = "0.24 (-2.25;  2.74)"
$ws.Range("G9").Value = "0.21 (-1.98;  2.40)"
$ws.Range("H9").Value = "clomipramina"
$ws.Range("I9").Value = "0.24 (-1.22;  1.71)"

# Linha 10
$ws.Range("A10").Value = "2.02 ( 1.32;  2.73)"
$ws.Range("B10").Value = "2.02 ( 1.06;  2.98)"
$ws.Range("C10").Value = "1.47 ( 0.27;  2.66)"
$ws.Range("D10").Value = "1.23 ( 0.31;  2.15)"
$ws.Range("E10").Value = "1.17 (-0.22;  2.56)"
$ws.Range("F10").Value = "0.39 (-1.67;  2.45)"
$ws.Range("G10").Value = "0.35 (-1.33;  2.04)"
$ws.Range("H10").Value = "0.15 (-1.26;  1.55)"
$ws.Range("J10").Value = "0.36 (-1.33;  2.04)"

# Linha 11
$ws.Range("A11").Value = "2.38 ( 0.55;  4.20)"
$ws.Range("B11").Value = "2.38 ( 0.44;  4.31)"
$ws.Range("C11").Value = "1.82 (-0.24;  3.89)"
$ws.Range("D11").Value = "1.59 (-0.33;  3.51)"
$ws.Range("E11").Value = "1.52 (-0.66;  3.71)"
$ws.Range("F11").Value = "0.75 (-1.92;  3.41)"
$ws.Range("G11").Value = "0.71 (-1.67;  3.09)"
$ws.Range("H11").Value = "0.50 (-1.69;  2.69)"
$ws.Range("I11").Value = "0.36 (-1.33;  2.04)"
